# customerTransactions: fix row 2, populate row 3, seed Order # for row 4
# (commit: "get revenue woks for manager class")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Order # 2): trim name, zero-out pepperoni count, add 3 cheese
# pizzas, and refresh the order total.
$ws.Range("B2").Value = "e"
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 0
$ws.Range("H2").Value = 48.26

# Row 3 (Order # 3) was a blank row stub; fill in the customer + order
# details. B3's "55" must stay TEXT (not auto-coerced to the number 55),
# matching the original data's inline-string cell. Assigning it directly
# would have Excel infer a number, so stage it in a text-formatted scratch
# cell and bring over only the *value* via PasteSpecial (values only),
# which leaves B3's own formatting untouched.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"
$helper.Value = "55"
$helper.Copy()
$ws.Range("B3").PasteSpecial(-4163)   # xlPasteValues
$helper.Clear()
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "cf"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 34
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 583.4400000000001

# Row 4 only gets its Order # populated.
$ws.Range("A4").Value = 4
